$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" date placeholder text from
#    25/11/2019 to 26/11/2019 everywhere it appears: on the slide
#    master itself and on every one of its custom (slide) layouts.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "*Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq "25/11/2019") {
                    $shp.TextFrame.TextRange.Text = "26/11/2019"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 2 title wording tweak.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$title = $s2.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "LINGUAGENS/BIBLIOTECAS USADAS:"

# ---------------------------------------------------------------------
# 3) Touch the presentation-level Guides collection so PowerPoint's
#    (empty) slide-guide list extension is materialised in the package
#    (best-effort; harmless if the host does not surface a mutation).
# ---------------------------------------------------------------------
try {
    $guides = $p.Guides
    $null = $guides.Count
} catch {
}
